$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Complete row 10 with PriceChange and UpDown values
$ws.Range("X10").Value = 0.6499990000000011
$ws.Range("Y10").Value = "Up"

# Widen column C to fit new "Strong Buy" verdict text
$ws.Columns.Item(3).ColumnWidth = 8.3

# Add new row 11 (traded row)
$ws.Range("A11").Value = 42654.883275462962
$ws.Range("B11").Value = 22
$ws.Range("C11").Value = "Strong Buy"
$ws.Range("D11").Value = 20
$ws.Range("E11").Value = 12138
$ws.Range("F11").Value = 1242
$ws.Range("G11").Value = 61
$ws.Range("H11").Value = 36
$ws.Range("I11").Value = 87
$ws.Range("J11").Value = 12
$ws.Range("K11").Value = 10936
$ws.Range("L11").Value = 199
$ws.Range("M11").Value = 117
$ws.Range("N11").Value = 36
$ws.Range("O11").Value = 5
$ws.Range("P11").Value = "Bag"
$ws.Range("Q11").Value = 58.438771163779279
$ws.Range("R11").Value = 0.49
$ws.Range("S11").Value = 0.0933
$ws.Range("T11").Value = 0.0249
$ws.Range("U11").Value = 2.34
$ws.Range("V11").Value = "N/A"
$ws.Range("W11").Value = 2

# Match number formats used by the other rows for these percentage columns
$ws.Range("S11").NumberFormat = "0.00%"
$ws.Range("T11").NumberFormat = "0.00%"
